$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 120.71429
$ws.Range("I5").Value = 120.71429
$ws.Range("K5").Value = 120.71429
$ws.Range("M5").Value = -5.714290000000005
$ws.Range("H32").Value = 2979.4
$ws.Range("I32").Value = 4997
$ws.Range("K32").Value = 4997
$ws.Range("M32").Value = -4671
$ws.Range("H43").Value = 5232.5
$ws.Range("J43").Value = 5879.2
$ws.Range("L43").Value = 5879.2
$ws.Range("N43").Value = -6017.2
$ws.Range("H51").Value = 7633.3335
$ws.Range("J51").Value = 7633.3335
$ws.Range("L51").Value = 7633.3335
$ws.Range("N51").Value = -8601.333500000001
$ws.Range("H64").Value = 4870
$ws.Range("J64").Value = 3175
$ws.Range("L64").Value = 3175
$ws.Range("N64").Value = -3671
$ws.Range("H67").Value = 4870
$ws.Range("J67").Value = 3175
$ws.Range("L67").Value = 3175
$ws.Range("N67").Value = -4891
$ws.Range("H74").Value = 3998.25
$ws.Range("I74").Value = 3998.25
$ws.Range("K74").Value = 3998.25
$ws.Range("M74").Value = -3062.25
$ws.Range("H77").Value = 3998.25
$ws.Range("I77").Value = 3998.25
$ws.Range("K77").Value = 19991.25
$ws.Range("M77").Value = -15311.25
$ws.Range("H100").Value = 1481
$ws.Range("I100").Value = 904.2857
$ws.Range("K100").Value = 904.2857
$ws.Range("M100").Value = -363.2857
$ws.Range("H113").Value = 3365.6667
$ws.Range("I113").Value = 3265
$ws.Range("J113").Value = 3466.3333
$ws.Range("K113").Value = 3265
$ws.Range("L113").Value = 3466.3333
$ws.Range("M113").Value = -11
$ws.Range("N113").Value = -9974.3333
$ws.Range("H132").Value = 3611.75
$ws.Range("I132").Value = 3611.75
$ws.Range("K132").Value = 10835.25
$ws.Range("M132").Value = -8305.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5465.533
$ws.Range("J63").Value = 6987.1113
$ws.Range("L63").Value = 6987.1113
$ws.Range("N63").Value = -8359.1113
$ws.Range("H66").Value = 5465.533
$ws.Range("J66").Value = 6987.1113
$ws.Range("L66").Value = 34935.5565
$ws.Range("N66").Value = -41799.5565
$ws.Range("H102").Value = 1388.25
$ws.Range("I102").Value = 1321.1428
$ws.Range("K102").Value = 1321.1428
$ws.Range("M102").Value = 300.8571999999999
$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 1993.3334
$ws.Range("I86").Value = 2682
$ws.Range("K86").Value = 2682
$ws.Range("M86").Value = -1559
$ws.Range("H89").Value = 1993.3334
$ws.Range("I89").Value = 2682
$ws.Range("K89").Value = 13410
$ws.Range("M89").Value = -7794
$ws.Range("H94").Value = 986.25
$ws.Range("I94").Value = 1019.0909
$ws.Range("K94").Value = 1019.0909
$ws.Range("M94").Value = -568.0909
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 355.81818
$ws.Range("I22").Value = 173.28572
$ws.Range("J22").Value = 675.25
$ws.Range("K22").Value = 173.28572
$ws.Range("L22").Value = 675.25
$ws.Range("M22").Value = 176.71428
$ws.Range("N22").Value = -1375.25
$ws.Range("H31").Value = 4745.8096
$ws.Range("I31").Value = 2829.6924
$ws.Range("K31").Value = 2829.6924
$ws.Range("M31").Value = -2534.6924
$ws.Range("H34").Value = 4745.8096
$ws.Range("I34").Value = 2829.6924
$ws.Range("K34").Value = 2829.6924
$ws.Range("M34").Value = -2627.6924
$ws.Range("H86").Value = 8415.916999999999
$ws.Range("I86").Value = 7170.2856
$ws.Range("K86").Value = 7170.2856
$ws.Range("M86").Value = -6047.2856
$ws.Range("H89").Value = 8415.916999999999
$ws.Range("I89").Value = 7170.2856
$ws.Range("K89").Value = 35851.428
$ws.Range("M89").Value = -30235.428
$ws.Range("H105").Value = 708.15
$ws.Range("I105").Value = 492.5
$ws.Range("J105").Value = 2649
$ws.Range("K105").Value = 492.5
$ws.Range("L105").Value = 2649
$ws.Range("M105").Value = 1254.5
$ws.Range("N105").Value = -6143
$ws.Range("H132").Value = 4501.3335
$ws.Range("I132").Value = 3998
$ws.Range("K132").Value = 11994
$ws.Range("M132").Value = -9464
$ws.Range("H134").Value = 3148.1538
$ws.Range("I134").Value = 2304.75
$ws.Range("K134").Value = 6914.25
$ws.Range("M134").Value = -4379.25
$ws.Range("H135").Value = 24999
$ws.Range("I135").Value = 24999
$ws.Range("K135").Value = 24999
$ws.Range("M135").Value = -19929
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 76967
$ws.Range("I2").Value = 125031.125
$ws.Range("K2").Value = 750186.75
$ws.Range("M2").Value = -750073.75
$ws.Range("H3").Value = 1478
$ws.Range("I3").Value = 1478
$ws.Range("K3").Value = 4434
$ws.Range("M3").Value = -4322
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4316.6924
$ws.Range("I80").Value = 3828.875
$ws.Range("J80").Value = 5097.2
$ws.Range("K80").Value = 3828.875
$ws.Range("L80").Value = 5097.2
$ws.Range("M80").Value = -2830.875
$ws.Range("N80").Value = -7093.2
$ws.Range("H83").Value = 4316.6924
$ws.Range("I83").Value = 3828.875
$ws.Range("J83").Value = 5097.2
$ws.Range("K83").Value = 19144.375
$ws.Range("L83").Value = 25486
$ws.Range("M83").Value = -14152.375
$ws.Range("N83").Value = -35470
$ws.Range("H92").Value = 11025
$ws.Range("J92").Value = 9472.223
$ws.Range("L92").Value = 9472.223
$ws.Range("N92").Value = -13216.223
$ws.Range("H95").Value = 1010000
$ws.Range("I95").Value = 2000000
$ws.Range("J95").Value = 20000
$ws.Range("K95").Value = 2000000
$ws.Range("L95").Value = 20000
$ws.Range("M95").Value = -1997254
$ws.Range("N95").Value = -25492
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1823.5
$ws.Range("I40").Value = 1823.5
$ws.Range("K40").Value = 1823.5
$ws.Range("M40").Value = -1687.5
$ws.Range("H61").Value = 1521.2222
$ws.Range("I61").Value = 1216.4
$ws.Range("K61").Value = 1216.4
$ws.Range("M61").Value = -1014.4
$ws.Range("H82").Value = 3650.5789
$ws.Range("I82").Value = 3670.5454
$ws.Range("K82").Value = 3670.5454
$ws.Range("M82").Value = -3309.5454
$ws.Range("H85").Value = 3650.5789
$ws.Range("I85").Value = 3670.5454
$ws.Range("K85").Value = 3670.5454
$ws.Range("M85").Value = -2422.5454
$ws.Range("H93").Value = 3045.1
$ws.Range("I93").Value = 3045.1
$ws.Range("K93").Value = 3045.1
$ws.Range("M93").Value = -1797.1
$ws.Range("H103").Value = 15000
$ws.Range("J103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -17344
$ws.Range("H113").Value = 1521.2222
$ws.Range("I113").Value = 1216.4
$ws.Range("K113").Value = 1216.4
$ws.Range("M113").Value = 953.5999999999999
$ws.Range("H132").Value = 4985.6665
$ws.Range("I132").Value = 4985.6665
$ws.Range("K132").Value = 14956.9995
$ws.Range("M132").Value = -12426.9995
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H45").Value = 20151.334
$ws.Range("J45").Value = 18637.572
$ws.Range("L45").Value = 18637.572
$ws.Range("N45").Value = -19619.572
$ws.Range("H62").Value = 7642.4287
$ws.Range("I62").Value = 5499.75
$ws.Range("K62").Value = 5499.75
$ws.Range("M62").Value = -4875.75
$ws.Range("H65").Value = 7642.4287
$ws.Range("I65").Value = 5499.75
$ws.Range("K65").Value = 27498.75
$ws.Range("M65").Value = -24378.75
$ws.Range("H74").Value = 35899
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 43848.5
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 43848.5
$ws.Range("M74").Value = -19064
$ws.Range("N74").Value = -45720.5
$ws.Range("H77").Value = 35899
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 43848.5
$ws.Range("K77").Value = 60000
$ws.Range("L77").Value = 131545.5
$ws.Range("M77").Value = -55320
$ws.Range("N77").Value = -140905.5
$ws.Range("H81").Value = 2311.0667
$ws.Range("I81").Value = 1853.75
$ws.Range("J81").Value = 4140.3335
$ws.Range("K81").Value = 3707.5
$ws.Range("L81").Value = 8280.666999999999
$ws.Range("M81").Value = -2646.5
$ws.Range("N81").Value = -10402.667
$ws.Range("H84").Value = 2311.0667
$ws.Range("I84").Value = 1853.75
$ws.Range("J84").Value = 4140.3335
$ws.Range("K84").Value = 18537.5
$ws.Range("L84").Value = 41403.335
$ws.Range("M84").Value = -13233.5
$ws.Range("N84").Value = -52011.335
